$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry row: date, hours spent, description
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 43256
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "Runnable thread erstellt "

$ws.Range("C11").Select()
